$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "het robotjs gebruikt" -> "robotjs gebruikt"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "het robotjs gebruikt",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "robotjs gebruikt", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Split the paragraph right after "kunnen installeren." into:
#      - end of current paragraph
#      - a new, empty paragraph
#      - a new paragraph that used to be the " Met de hulp..." run
# ---------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("kunnen installeren.") | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertParagraphAfter() | Out-Null

$rng.Collapse(0) | Out-Null
$rng.MoveEnd(1, 1) | Out-Null
$rng.Collapse(0) | Out-Null
$rng.InsertAfter([char]13) | Out-Null

# ---------------------------------------------------------------------
# 3) Rewrite the old run's content (now living in its own paragraph) into
#    the new text: "het computer os" -> "de computers OS" and append the
#    extra clause before the final period.
# ---------------------------------------------------------------------
$apos = [char]0x2019

$rng2 = $d.Content
$rng2.Find.Execute(" Met de hulp van mijn leraar Berend ben ik er ook achter gekomen waarom games de input niet konden lezen, dit kwam omdat ze naar keyDown events luisteren, maar node-key-sender verstuurd een keyUp event naar het computer os. Ook heb gevonden dat de emulator die ik wou gebruiken voor de test alleen luistert naar een geregistreerde controller(voorbeeld: Xbox, PlayStation en keyboard) dat betekent dat geen van deze programma" + $apos + "s input geregistreerd wordt.") | Out-Null

$rng2.Text = "Met de hulp van mijn leraar Berend ben ik er ook achter gekomen waarom games de input niet konden lezen, dit kwam omdat ze naar keyDown events luisteren, maar node-key-sender verstuurd een keyUp event naar de computers OS. Ook heb gevonden dat de emulator die ik wou gebruiken voor de test alleen luistert naar een geregistreerde controller(voorbeeld: Xbox, PlayStation en keyboard) dat betekent dat geen van deze programma" + $apos + "s input geregistreerd wordt omdat het niet van deze geregistreerd programma" + $apos + "s kwam."

Write-Output "done"
